# Brewery Parts.xlsx edit script
# Commit message: "Renamed Lable Tun.lsl Generalise Tun Contents for MT + Kettle"
#
# This applies the cell value / formatting / layout changes described by the
# authoritative xml diff against xl/worksheets/sheet1.xml, xl/styles.xml and
# xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column width tweaks (col C widened, col E widened slightly)
# ---------------------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 14.02
$ws.Columns(5).ColumnWidth = 20.88

# ---------------------------------------------------------------------------
# Row 2 : Base -- "!Animate Brewery" -> "Brewery"
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Brewery"

# ---------------------------------------------------------------------------
# Row 4 : Mill -- "Animate Mill" -> "Grain Mill"
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "Grain Mill"

# ---------------------------------------------------------------------------
# Row 5 : HLT -- new Desc "Hot Liquor Tank", "Label Tun" -> "Label Object from Desc"
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "Hot Liquor Tank"
$ws.Range("E5").Value = "Label Object from Desc"

# ---------------------------------------------------------------------------
# Row 6 : HLT Lid -- "Animate HLT Lid" -> "Tun Lid", drop the old Notes comment
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = "Tun Lid"
$ws.Range("I6").ClearContents()

# ---------------------------------------------------------------------------
# Row 7 : HLT contents -- new Script/s "HLT Steam"
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "HLT Steam"

# ---------------------------------------------------------------------------
# Row 9 : MT -- Desc "Mash Tun" -> "Brewery", "Label Tun" -> "Label Object from Desc"
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "Brewery"
$ws.Range("E9").Value = "Label Object from Desc"

# ---------------------------------------------------------------------------
# Row 10 : MT Lid -- "Animate MT Lid" -> "Tun Lid", drop the old Notes comment
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = "Tun Lid"
$ws.Range("I10").ClearContents()

# ---------------------------------------------------------------------------
# Row 11 : MT Contents -- "Animate MT Contents" -> "Tun Contents", new F/I notes
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "Tun Contents"
$ws.Range("F11").Value = "Grain/Mash ??"
$ws.Range("I11").Value = "Differentiate grain to mash?"

# ---------------------------------------------------------------------------
# Row 13 : Kettle -- "Label Tun" -> "Label Object from Desc"
# ---------------------------------------------------------------------------
$ws.Range("E13").Value = "Label Object from Desc"

# ---------------------------------------------------------------------------
# Row 14 : Kettle lid -- "Animate Kettle Lid" -> "Tun Lid", drop old Notes comment
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = "Tun Lid"
$ws.Range("I14").ClearContents()

# ---------------------------------------------------------------------------
# Row 15 : Kettle contents -- new Script/s "Tun Contents"
# ---------------------------------------------------------------------------
$ws.Range("E15").Value = "Tun Contents"

# ---------------------------------------------------------------------------
# Row 16 : Chimney -- Name cell renamed to "Mash Tun", new Notes "add puffer to top"
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "Mash Tun"
$ws.Range("I16").Value = "add puffer to top"

# ---------------------------------------------------------------------------
# Row 17 : Pipe1 -- new Script/s "Pump Assembly 1"
# ---------------------------------------------------------------------------
$ws.Range("E17").Value = "Pump Assembly 1"

# ---------------------------------------------------------------------------
# Row 18 : Pump1 -- new wrapped Script/s, taller row
# ---------------------------------------------------------------------------
$ws.Range("E18").Value = "Pump Assembly 1`nLabel Obect from Desc"
$ws.Range("E18").WrapText = $true
$ws.Rows(18).RowHeight = 31.5

# ---------------------------------------------------------------------------
# Row 19 : Pipe2a -- new Script/s "Pump Assembly 2"
# ---------------------------------------------------------------------------
$ws.Range("E19").Value = "Pump Assembly 2"

# ---------------------------------------------------------------------------
# Row 20 : pump2 -- new wrapped Script/s, taller row, Notes highlighted red
# ---------------------------------------------------------------------------
$ws.Range("E20").Value = "Pump Assembly 2`nLabel Object from Desc"
$ws.Range("E20").WrapText = $true
$ws.Rows(20).RowHeight = 31.5
$ws.Range("I20").Value = "How do we reverse the flow?"
$ws.Range("I20").Interior.Color = 255

# ---------------------------------------------------------------------------
# Row 21 : pipe2b -- new Script/s "Pump Assembly 2", Notes highlighted red
# ---------------------------------------------------------------------------
$ws.Range("E21").Value = "Pump Assembly 2"
$ws.Range("I21").Value = "How do we reverse the flow?"
$ws.Range("I21").Interior.Color = 255

# ---------------------------------------------------------------------------
# Row 22 : pipe3a -- new Script/s "Chiller circuit"
# ---------------------------------------------------------------------------
$ws.Range("E22").Value = "Chiller circuit"

# ---------------------------------------------------------------------------
# Row 23 : Chiller -- new wrapped Script/s, shorter row
# ---------------------------------------------------------------------------
$ws.Range("E23").Value = "Chiller circuit`nLabel Object from Desc"
$ws.Range("E23").WrapText = $true
$ws.Rows(23).RowHeight = 29.25

# ---------------------------------------------------------------------------
# Row 24 : pipe3b -- new Script/s "Chiller circuit"
# ---------------------------------------------------------------------------
$ws.Range("E24").Value = "Chiller circuit"

# ---------------------------------------------------------------------------
# Final selection, as left by the author
# ---------------------------------------------------------------------------
$ws.Range("F11").Select()
